$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "3>2>1" priority output values (column U)
$ws.Range("U6").Value = 24
$ws.Range("U7").Value = 74
$ws.Range("U9").Value = 1

# Move the active selection to U10 (single cell)
$ws.Range("U10").Select()
